$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.884.79"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.85"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.71"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5087"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07177"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8897"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.56"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.50"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07490"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.44"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008482"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.15"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.937.94"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.003"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.115.07"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.33"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.366"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.43"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.779"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.092"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.55"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.688"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.697"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09104"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05027"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.965"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.152"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.209"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.504"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5547"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.073"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.579"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.85"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.579"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1483"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4744"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.10"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.10"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.555"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.95"
$ws.Range("E51").Value = "  -1.23%  "
